# Criada a classe gráfico
# Update the statistics row (row 18) of the homicidiosNegrosPaisQTD sheet:
#  - D18: AVEDEV(...) -> MEDIAN(...)
#  - E18: formula removed, replaced by literal 0
#  - F18: now holds the STDEV.P formula (previously in E18)
#  - G18: now holds the VAR.P formula (previously in F18); the old
#         coefficient-of-variation formula (E18/C18) that used to live
#         in G18 is gone.
#  - H18 (VAR.S) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Formula = "=MEDIAN(B1:B18)"
$ws.Range("E18").Value = 0
$ws.Range("F18").Formula = "=_xlfn.STDEV.P(B1:B18)"
$ws.Range("G18").Formula = "=_xlfn.VAR.P(B1:B18)"

# Reflect the final selection left by the author (cell E19, just below
# the edited row).
$ws.Range("E19").Select()
